# Add the "student ID" rectangle shape to slide 1.
#
# Slides 2 and 3 already contain an identical shape (a thin rectangle with an
# accent1-colored outline, positioned top-right, containing the text
# "IT14098888"). Slide 1 is missing it, so we reproduce it there by copying
# the existing shape from slide 2 and pasting it onto slide 1 - this gives us
# a faithful clone (fill/line/paragraph/run formatting, bodyPr, lstStyle,
# endParaRPr, etc.) instead of hand-building the shape from scratch.

$p = $ppt.ActivePresentation
$sourceSlide = $p.Slides.Item(2)
$targetSlide = $p.Slides.Item(1)

# "Rectangle 3" is the IT14098888 badge shape already present on slide 2.
$template = $sourceSlide.Shapes.Item("Rectangle 3")
$template.Copy()

$newShape = $targetSlide.Shapes.Paste().Item(1)
$newShape.Name = "Rectangle 6"
